$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.2039992849870629
$ws.Range("D2").Value = 0.03013948496803209
$ws.Range("E2").Value = 0.2041107199060548
$ws.Range("F2").Value = 0.8009902855880071
$ws.Range("G2").Value = 0.002435002248668664
$ws.Range("L2").Value = 0.4352601544044887
$ws.Range("M2").Value = 10.94388192429562
$ws.Range("O2").Value = 2.777273334386024

$ws.Range("C3").Value = 0.2133631626660701
$ws.Range("D3").Value = 0.03027348549026243
$ws.Range("E3").Value = 0.1898001546332324
$ws.Range("F3").Value = 0.8205966484033524
$ws.Range("G3").Value = 0.002441226194074987
$ws.Range("L3").Value = 0.3874367492097406
$ws.Range("M3").Value = 9.62274159506353
$ws.Range("O3").Value = 2.871687663754955

$ws.Range("C4").Value = 0.2195782688237529
$ws.Range("D4").Value = 0.03040146940685418
$ws.Range("E4").Value = 0.1812430611963762
$ws.Range("F4").Value = 0.8342316218690158
$ws.Range("G4").Value = 0.00244520908604774
$ws.Range("L4").Value = 0.3582604089315282
$ws.Range("M4").Value = 8.808708495625979
$ws.Range("O4").Value = 2.935545167186689

$ws.Range("C5").Value = 0.2222258444194836
$ws.Range("D5").Value = 0.03046488571516193
$ws.Range("E5").Value = 0.177811636452752
$ws.Range("F5").Value = 0.8401823622189042
$ws.Range("G5").Value = 0.002446872963367254
$ws.Range("L5").Value = 0.3464147563801134
$ws.Range("M5").Value = 8.476230920908563
$ws.Range("O5").Value = 2.963021456590127

$ws.Range("C6").Value = 0.222672346454786
$ws.Range("D6").Value = 0.03047608957329118
$ws.Range("E6").Value = 0.1772451542262772
$ws.Range("F6").Value = 0.8411940915783092
$ws.Range("G6").Value = 0.00244715172101112
$ws.Range("L6").Value = 0.3444503609382537
$ws.Range("M6").Value = 8.420976854278535
$ws.Range("O6").Value = 2.967670914928107

$ws.Range("C7").Value = 0.2196135127385332
$ws.Range("D7").Value = 0.03040227935974826
$ws.Range("E7").Value = 0.1811965610982469
$ws.Range("F7").Value = 0.8343102881222748
$ws.Range("G7").Value = 0.002445231360027623
$ws.Range("L7").Value = 0.3581004804090639
$ws.Range("M7").Value = 8.80422767821608
$ws.Range("O7").Value = 2.935909870153353

$ws.Range("C8").Value = 0.2071299489273457
$ws.Range("D8").Value = 0.03017605889399277
$ws.Range("E8").Value = 0.1991274710974551
$ws.Range("F8").Value = 0.8074148935545153
$ws.Range("G8").Value = 0.002437114921382603
$ws.Range("L8").Value = 0.4187297282005602
$ws.Range("M8").Value = 10.48891995815609
$ws.Range("O8").Value = 2.808589618100143

$ws.Range("C9").Value = 0.1864471037822852
$ws.Range("D9").Value = 0.0301061488650447
$ws.Range("E9").Value = 0.2362147887839399
$ws.Range("F9").Value = 0.7676635705779802
$ws.Range("G9").Value = 0.002422467520876508
$ws.Range("L9").Value = 0.5392741153589782
$ws.Range("M9").Value = 13.77205026641826
$ws.Range("O9").Value = 2.606821435530009

$ws.Range("C10").Value = 0.1737126929135542
$ws.Range("D10").Value = 0.03029882609312295
$ws.Range("E10").Value = 0.2647913531924075
$ws.Range("F10").Value = 0.7468281424850289
$ws.Range("G10").Value = 0.002412462915595785
$ws.Range("L10").Value = 0.6291002719380572
$ws.Range("M10").Value = 16.17507416934058
$ws.Range("O10").Value = 2.489459369112552

$ws.Range("C11").Value = 0.1684873252146843
$ws.Range("D11").Value = 0.03044327573518046
$ws.Range("E11").Value = 0.2781170592833462
$ws.Range("F11").Value = 0.7392669193691432
$ws.Range("G11").Value = 0.00240807211309578
$ws.Range("L11").Value = 0.6703011496467184
$ws.Range("M11").Value = 17.26716896728499
$ws.Range("O11").Value = 2.443143952223352

$ws.Range("C12").Value = 0.166593151532922
$ws.Range("D12").Value = 0.03050647773449811
$ws.Range("E12").Value = 0.2832132828241214
$ws.Range("F12").Value = 0.7366878144456166
$ws.Range("G12").Value = 0.002406432177266993
$ws.Range("L12").Value = 0.6859569797923086
$ws.Range("M12").Value = 17.68064036012788
$ws.Range("O12").Value = 2.426654246453353

$ws.Range("C13").Value = 0.1669972867256462
$ws.Range("D13").Value = 0.03049248239945257
$ws.Range("E13").Value = 0.2821134403621528
$ws.Range("F13").Value = 0.7372304945221941
$ws.Range("G13").Value = 0.002406784359039974
$ws.Range("L13").Value = 0.6825827288999733
$ws.Range("M13").Value = 17.59159432188426
$ws.Range("O13").Value = 2.430158432137517

$ws.Range("C14").Value = 0.1683297784115894
$ws.Range("D14").Value = 0.03044830302690116
$ws.Range("E14").Value = 0.2785353063080436
$ws.Range("F14").Value = 0.7390489901932398
$ws.Range("G14").Value = 0.002407936739853303
$ws.Range("L14").Value = 0.6715880513109767
$ws.Range("M14").Value = 17.30118670957518
$ws.Range("O14").Value = 2.44176612671032

$ws.Range("C15").Value = 0.1691570696538349
$ws.Range("D15").Value = 0.03042235934307769
$ws.Range("E15").Value = 0.276350213944923
$ws.Range("F15").Value = 0.7402001345816984
$ws.Range("G15").Value = 0.002408645562559739
$ws.Range("L15").Value = 0.6648606832648056
$ws.Range("M15").Value = 17.12329528124087
$ws.Range("O15").Value = 2.449013736126034

$ws.Range("C16").Value = 0.1740658528458852
$ws.Range("D16").Value = 0.03029055708368844
$ws.Range("E16").Value = 0.2639273156363799
$ws.Range("F16").Value = 0.7473615825069686
$ws.Range("G16").Value = 0.002412753066694756
$ws.Range("L16").Value = 0.6264149512862218
$ws.Range("M16").Value = 16.10368787223479
$ws.Range("O16").Value = 2.49263122242067

$ws.Range("C17").Value = 0.1772246579312196
$ws.Range("D17").Value = 0.03022448052612248
$ws.Range("E17").Value = 0.2563919853904082
$ws.Range("F17").Value = 0.7522518588908014
$ws.Range("G17").Value = 0.00241531375131959
$ws.Range("L17").Value = 0.6029200694904944
$ws.Range("M17").Value = 15.47796999089621
$ws.Range("O17").Value = 2.521223131531883

$ws.Range("C18").Value = 0.1790947590317273
$ws.Range("D18").Value = 0.03019180019878576
$ws.Range("E18").Value = 0.2520883530783209
$ws.Range("F18").Value = 0.7552445240276882
$ws.Range("G18").Value = 0.002416801699614188
$ws.Range("L18").Value = 0.5894379575112509
$ws.Range("M18").Value = 15.11796992133003
$ws.Range("O18").Value = 2.538332250869217

$ws.Range("C19").Value = 0.1797369994963702
$ws.Range("D19").Value = 0.03018163991652756
$ws.Range("E19").Value = 0.2506363603193194
$ws.Range("F19").Value = 0.7562884323311607
$ws.Range("G19").Value = 0.002417308097231744
$ws.Range("L19").Value = 0.5848784057665455
$ws.Range("M19").Value = 14.99606024114087
$ws.Range("O19").Value = 2.544238158334679

$ws.Range("C20").Value = 0.1768828666826607
$ws.Range("D20").Value = 0.03023096099195755
$ws.Range("E20").Value = 0.2571909561466441
$ws.Range("F20").Value = 0.7517125998855789
$ws.Range("G20").Value = 0.002415039600271071
$ws.Range("L20").Value = 0.60541784033623
$ws.Range("M20").Value = 15.54458897165904
$ws.Range("O20").Value = 2.518110544372632

$ws.Range("C21").Value = 0.1679360743027445
$ws.Range("D21").Value = 0.03046104607875577
$ws.Range("E21").Value = 0.2795849058588686
$ws.Range("F21").Value = 0.7385070727795551
$ws.Range("G21").Value = 0.002407597642144289
$ws.Range("L21").Value = 0.6748159463914192
$ws.Range("M21").Value = 17.38648807711735
$ws.Range("O21").Value = 2.438327939754686

$ws.Range("C22").Value = 0.1625831627753627
$ws.Range("D22").Value = 0.03066114509451268
$ws.Range("E22").Value = 0.294514100796448
$ws.Range("F22").Value = 0.7315369476222884
$ws.Range("G22").Value = 0.002402866447315217
$ws.Range("L22").Value = 0.720488826590298
$ws.Range("M22").Value = 18.58983955332758
$ws.Range("O22").Value = 2.392314787228486

$ws.Range("C23").Value = 0.1653938757212501
$ws.Range("D23").Value = 0.03054968327973739
$ws.Range("E23").Value = 0.2865181545308957
$ws.Range("F23").Value = 0.7351022316056088
$ws.Range("G23").Value = 0.002405379545398211
$ws.Range("L23").Value = 0.6960814799967352
$ws.Range("M23").Value = 17.9476032254251
$ws.Range("O23").Value = 2.416301296379231

$ws.Range("C24").Value = 0.177037222338555
$ws.Range("D24").Value = 0.03022801466502045
$ws.Range("E24").Value = 0.2568296528429101
$ws.Range("F24").Value = 0.7519558351881983
$ws.Range("G24").Value = 0.002415163494794522
$ws.Range("L24").Value = 0.604288519389371
$ws.Range("M24").Value = 15.51447136584699
$ws.Range("O24").Value = 2.519515656632478

$ws.Range("C25").Value = 0.1916210728167584
$ws.Range("D25").Value = 0.03008343431472582
$ws.Range("E25").Value = 0.2259602292246967
$ws.Range("F25").Value = 0.7769806611570615
$ws.Range("G25").Value = 0.002426295880433836
$ws.Range("L25").Value = 0.5064617793764512
$ws.Range("M25").Value = 12.88574302349832
$ws.Range("O25").Value = 2.426654246453353

